$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rename the "Full scan" caption for D6 - it is now about file selection instead of upload.
$ws.Range("D6").Value = "Full scan complete and file selected?"

# Insert three new rows right before the old "parameters" block (old row 23) to hold
# the new "peak_shape" entries for carbon, hydrogen and nitrogen.
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(23).Insert()

# Row 23: peak_shape for carbon
$ws.Range("A23").Value = "sensitivity"
$ws.Range("B23").Value = "carbon"
$ws.Range("C23").Value = "peak_shape"
$ws.Range("D23").Value = "Peak shape scan complete & selected?"
$ws.Range("E23").Value = "bool"
$ws.Range("G23").Value = "Run a peak shape scan (from ? To ?), split out, BF on"

# Row 24: peak_shape for hydrogen
$ws.Range("A24").Value = "sensitivity"
$ws.Range("B24").Value = "hydrogen"
$ws.Range("C24").Value = "peak_shape"
$ws.Range("D24").Value = "Peak shape scan complete & selected?"
$ws.Range("E24").Value = "bool"
$ws.Range("G24").Value = "Run a peak shape scan (from ? To ?), split out, BF on"

# Row 25: peak_shape for nitrogen
$ws.Range("A25").Value = "sensitivity"
$ws.Range("B25").Value = "nitrogen"
$ws.Range("C25").Value = "peak_shape"
$ws.Range("D25").Value = "Peak shape scan complete & selected?"
$ws.Range("E25").Value = "bool"
$ws.Range("G25").Value = "Run a peak shape scan (from ? To ?), split out, BF on"

# The Caption column now needs to be a bit wider to fit the new text.
$ws.Columns.Item(4).ColumnWidth = 31.83

# Match the author's final selection in the saved workbook.
$ws.Range("D25:G25").Select() | Out-Null
